# Edit script: add "metadata" worksheet after "data", and refresh the
# per-row "time_taken" timestamps on the "data" sheet (panelapp re-scrape).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Add the new "metadata" sheet right after "data"
# ---------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Differences of Sex Development"
$metaSheet.Range("C2").Value = 99

# data_version ("0.215") must stay a text value, not be coerced to a number.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.215"
$metaSheet.Range("D2").ClearFormats()

$metaSheet.Range("E2").Value = "2021-09-29T03:15:02.114871Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:43.780354"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/99/?format=json"

# Match the header styling used on "data" (bold / bordered / centered),
# and the numeric-row style used for column A, by copying the existing
# formats over (reuses the same style records, doesn't create new ones).
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Refresh "time_taken" (column F) on "data" for every gene row (2..105)
#    with the values captured from the re-run of the scrape.
# ---------------------------------------------------------------------
$dataTimeTaken = @(
  "2021-10-05 14:33:43.784085",
  "2021-10-05 14:33:43.784093",
  "2021-10-05 14:33:43.784097",
  "2021-10-05 14:33:43.784099",
  "2021-10-05 14:33:43.784102",
  "2021-10-05 14:33:43.784105",
  "2021-10-05 14:33:43.784108",
  "2021-10-05 14:33:43.784110",
  "2021-10-05 14:33:43.784113",
  "2021-10-05 14:33:43.784115",
  "2021-10-05 14:33:43.784118",
  "2021-10-05 14:33:43.784120",
  "2021-10-05 14:33:43.784123",
  "2021-10-05 14:33:43.784125",
  "2021-10-05 14:33:43.784128",
  "2021-10-05 14:33:43.784130",
  "2021-10-05 14:33:43.784133",
  "2021-10-05 14:33:43.784136",
  "2021-10-05 14:33:43.784138",
  "2021-10-05 14:33:43.784141",
  "2021-10-05 14:33:43.784143",
  "2021-10-05 14:33:43.784145",
  "2021-10-05 14:33:43.784148",
  "2021-10-05 14:33:43.784150",
  "2021-10-05 14:33:43.784153",
  "2021-10-05 14:33:43.784156",
  "2021-10-05 14:33:43.784158",
  "2021-10-05 14:33:43.784161",
  "2021-10-05 14:33:43.784163",
  "2021-10-05 14:33:43.784166",
  "2021-10-05 14:33:43.784168",
  "2021-10-05 14:33:43.784171",
  "2021-10-05 14:33:43.784174",
  "2021-10-05 14:33:43.784176",
  "2021-10-05 14:33:43.784179",
  "2021-10-05 14:33:43.784181",
  "2021-10-05 14:33:43.784183",
  "2021-10-05 14:33:43.784186",
  "2021-10-05 14:33:43.784189",
  "2021-10-05 14:33:43.784191",
  "2021-10-05 14:33:43.784195",
  "2021-10-05 14:33:43.784197",
  "2021-10-05 14:33:43.784200",
  "2021-10-05 14:33:43.784202",
  "2021-10-05 14:33:43.784205",
  "2021-10-05 14:33:43.784207",
  "2021-10-05 14:33:43.784210",
  "2021-10-05 14:33:43.784212",
  "2021-10-05 14:33:43.784215",
  "2021-10-05 14:33:43.784217",
  "2021-10-05 14:33:43.784220",
  "2021-10-05 14:33:43.784223",
  "2021-10-05 14:33:43.784226",
  "2021-10-05 14:33:43.784228",
  "2021-10-05 14:33:43.784231",
  "2021-10-05 14:33:43.784233",
  "2021-10-05 14:33:43.784236",
  "2021-10-05 14:33:43.784238",
  "2021-10-05 14:33:43.784241",
  "2021-10-05 14:33:43.784244",
  "2021-10-05 14:33:43.784246",
  "2021-10-05 14:33:43.784249",
  "2021-10-05 14:33:43.784251",
  "2021-10-05 14:33:43.784254",
  "2021-10-05 14:33:43.784257",
  "2021-10-05 14:33:43.784260",
  "2021-10-05 14:33:43.784263",
  "2021-10-05 14:33:43.784265",
  "2021-10-05 14:33:43.784268",
  "2021-10-05 14:33:43.784270",
  "2021-10-05 14:33:43.784273",
  "2021-10-05 14:33:43.784275",
  "2021-10-05 14:33:43.784278",
  "2021-10-05 14:33:43.784281",
  "2021-10-05 14:33:43.784283",
  "2021-10-05 14:33:43.784286",
  "2021-10-05 14:33:43.784291",
  "2021-10-05 14:33:43.784294",
  "2021-10-05 14:33:43.784297",
  "2021-10-05 14:33:43.784299",
  "2021-10-05 14:33:43.784302",
  "2021-10-05 14:33:43.784304",
  "2021-10-05 14:33:43.784307",
  "2021-10-05 14:33:43.784310",
  "2021-10-05 14:33:43.784312",
  "2021-10-05 14:33:43.784315",
  "2021-10-05 14:33:43.784317",
  "2021-10-05 14:33:43.784320",
  "2021-10-05 14:33:43.784323",
  "2021-10-05 14:33:43.784325",
  "2021-10-05 14:33:43.784328",
  "2021-10-05 14:33:43.784330",
  "2021-10-05 14:33:43.784334",
  "2021-10-05 14:33:43.784337",
  "2021-10-05 14:33:43.784340",
  "2021-10-05 14:33:43.784343",
  "2021-10-05 14:33:43.784345",
  "2021-10-05 14:33:43.784348",
  "2021-10-05 14:33:43.784351",
  "2021-10-05 14:33:43.784353",
  "2021-10-05 14:33:43.784356",
  "2021-10-05 14:33:43.784359",
  "2021-10-05 14:33:43.784361",
  "2021-10-05 14:33:43.784364"
)

for ($i = 0; $i -lt $dataTimeTaken.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $dataTimeTaken[$i]
}

$wb.Application.CutCopyMode = $false

# Keep "data" as the active/selected sheet (matches the unchanged activeTab
# in the workbook's bookViews).
$dataSheet.Activate()

Write-Output "metadata sheet added; time_taken refreshed for $($dataTimeTaken.Length) rows"
